# Update latest output (run 96)

$wb = $excel.ActiveWorkbook

# --- Schedule sheet ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 719.6379644999998
$wsSchedule.Range("F2").Value = 11.89877586805555

# --- Detailed sheet ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B8").Value = 36.2
$wsDetailed.Range("B9").Value = 48.37963
$wsDetailed.Range("B10").Value = 57.08

$wsDetailed.Range("B11").Value = 61.81795
$wsDetailed.Range("C11").Value = "historical"

$wsDetailed.Range("B12").Value = 60.49589
$wsDetailed.Range("C12").Value = "historical"

$wsDetailed.Range("B13").Value = 71.97167

$wsDetailed.Range("B16").Value = 36.06
$wsDetailed.Range("B17").Value = 8.797459999999999
$wsDetailed.Range("B18").Value = 0.7
$wsDetailed.Range("B19").Value = 0.7
$wsDetailed.Range("B20").Value = 0.7
$wsDetailed.Range("B21").Value = -2.83936

$wsDetailed.Range("B24").Value = -5.74313
$wsDetailed.Range("B25").Value = 0
$wsDetailed.Range("B26").Value = -2.83936
$wsDetailed.Range("B27").Value = -4.6323
$wsDetailed.Range("B28").Value = -4.73809
$wsDetailed.Range("B29").Value = -2.69309
$wsDetailed.Range("B30").Value = -0.93819
$wsDetailed.Range("B31").Value = 0.00003
$wsDetailed.Range("B32").Value = 0.51

$wsDetailed.Range("B34").Value = 1.67481
$wsDetailed.Range("B35").Value = -2.49115

$wsDetailed.Range("B37").Value = 0.009650000000000001
$wsDetailed.Range("B38").Value = 12.0924
$wsDetailed.Range("B39").Value = 42.36428
$wsDetailed.Range("B40").Value = 56.7552
$wsDetailed.Range("B41").Value = 59.19183
$wsDetailed.Range("B42").Value = 65

$wsDetailed.Range("B45").Value = 62.96274
$wsDetailed.Range("B46").Value = 61.15479
$wsDetailed.Range("B47").Value = 58.98416
$wsDetailed.Range("B48").Value = 58.49624
$wsDetailed.Range("B49").Value = 62.01821
